$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number need to be pre-formatted
# as Text so Excel stores them as strings (matching the original inlineStr
# cell type) instead of silently converting them to a numeric value.
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '30.198.92'
$ws.Range('E2').Value = '  -1.64%  '
$ws.Range('D3').Value = '1.847.86'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '232.81'
$ws.Range('E5').Value = '  -2.59%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '0.4707'
$ws.Range('E7').Value = '  -2.31%  '
$ws.Range('D8').Value = '0.2713'
$ws.Range('E8').Value = '  -4.62%  '
$ws.Range('D9').Value = '0.06370'
$ws.Range('E9').Value = '  -2.80%  '
$ws.Range('D10').Value = '1.859.91'
$ws.Range('E10').Value = '  -10.60%  '
$ws.Range('D11').Value = '0.07422'
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('D12').Value = '16.20'
$ws.Range('E12').Value = '  -2.96%  '
$ws.Range('D13').Value = '4.927'
$ws.Range('E13').Value = '  -3.66%  '
$ws.Range('D14').Value = '84.94'
$ws.Range('E14').Value = '  -3.57%  '
$ws.Range('D15').Value = '0.6258'
$ws.Range('E15').Value = '  -6.14%  '
$ws.Range('D16').Value = '30.152.84'
$ws.Range('E16').Value = '  -1.69%  '
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').Value = '227.78'
$ws.Range('E18').Value = '  -1.81%  '
$ws.Range('E19').Value = '  -5.48%  '
$ws.Range('D20').Value = '0.000007318'
$ws.Range('E20').Value = '  -3.85%  '
$ws.Range('D21').Value = '2.090.83'
$ws.Range('E21').Value = '  -3.63%  '
$ws.Range('D22').Value = '0.9995'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').Value = '4.917'
$ws.Range('E23').Value = '  -7.02%  '
$ws.Range('D24').Value = '5.906'
$ws.Range('E24').Value = '  -5.02%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '166.84'
$ws.Range('E25').Value = '  -1.71%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '9.197'
$ws.Range('E26').Value = '  -1.33%  '
$ws.Range('D27').Value = '17.85'
$ws.Range('E27').Value = '  -4.73%  '
$ws.Range('D28').Value = '1.866'
$ws.Range('E28').Value = '  -5.33%  '
$ws.Range('D29').Value = '0.1024'
$ws.Range('E29').Value = '  +0.90%  '
$ws.Range('D30').Value = '1.380'
$ws.Range('E30').Value = '  -1.70%  '
$ws.Range('D31').Value = '4.107'
$ws.Range('E31').Value = '  -5.92%  '
$ws.Range('E32').Value = '  -3.94%  '
$ws.Range('D33').Value = '0.04881'
$ws.Range('E33').Value = '  -4.41%  '
$ws.Range('D34').Value = '1.149'
$ws.Range('E34').Value = '  -5.88%  '
$ws.Range('D35').Value = '0.7063'
$ws.Range('E35').Value = '  -6.87%  '
$ws.Range('D36').Value = '0.9999'
$ws.Range('E36').Value = '  -0.62%  '
$ws.Range('D37').Value = '2.698'
$ws.Range('E37').Value = '  -0.66%  '
$ws.Range('D38').Value = '0.01847'
$ws.Range('E38').Value = '  -2.04%  '
$ws.Range('D39').Value = '2.630'
$ws.Range('E39').Value = '  -1.13%  '
$ws.Range('D40').Value = '0.9042'
$ws.Range('E40').Value = '  -1.95%  '
$ws.Range('D41').Value = '1.944'
$ws.Range('E41').Value = '  -6.75%  '
$ws.Range('D42').Value = '104.68'
$ws.Range('E42').Value = '  -2.64%  '
$ws.Range('D43').Value = '0.9980'
$ws.Range('E43').Value = '  -0.68%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '0.4074'
$ws.Range('E44').Value = '  -5.36%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '5.536'
$ws.Range('E45').Value = '  -3.87%  '
$ws.Range('D46').Value = '7.038'
$ws.Range('E46').Value = '  -5.39%  '
$ws.Range('D47').Value = '60.03'
$ws.Range('E47').Value = '  -7.33%  '
$ws.Range('D48').Value = '0.1189'
$ws.Range('E48').Value = '  -6.74%  '
$ws.Range('D49').Value = '8.627'
$ws.Range('E49').Value = '  -4.29%  '
$ws.Range('D50').Value = '33.08'
$ws.Range('E50').Value = '  -2.35%  '
$ws.Range('D51').Value = '1.383'
$ws.Range('E51').Value = '  -6.99%  '
